# Apply the workbook changes described by the commit:
# - Add a new "rémunération id" column (N on Chauffeur, I on Collab)
# - Renumber a couple of ids on Chauffeur and tweak a "service" column
# - Add two new rows of data on Collab

$wb = $excel.ActiveWorkbook
$wsChauffeur = $wb.Worksheets.Item("Chauffeur")
$wsCollab = $wb.Worksheets.Item("Collab")

# ---------------------------------------------------------------
# Chauffeur sheet (sheet1)
# ---------------------------------------------------------------

# New header for column N
$wsChauffeur.Range("N1").Value = "Id de la rémunération "

# Row 2: A 611 -> 621, C 24 -> 11, new N2 = 40
$wsChauffeur.Range("A2").Value = 621
$wsChauffeur.Range("C2").Value = 11
$wsChauffeur.Range("N2").Value = 40

# Row 3: A 607 -> 622, C 3 -> 11, new N3 = 41
$wsChauffeur.Range("A3").Value = 622
$wsChauffeur.Range("C3").Value = 11
$wsChauffeur.Range("N3").Value = 41

# Row 4: A 608 -> 623, C 3 -> 2, new N4 = 42
$wsChauffeur.Range("A4").Value = 623
$wsChauffeur.Range("C4").Value = 2
$wsChauffeur.Range("N4").Value = 42

# Widen the new column (closest achievable width to 20.7109375 chars)
$wsChauffeur.Columns.Item(14).ColumnWidth = 19.8

# Update the selection shown on the Chauffeur tab to D10:D11
$wsChauffeur.Range("D10:D11").Select()

# ---------------------------------------------------------------
# Collab sheet (sheet2)
# ---------------------------------------------------------------

# New header for column I
$wsCollab.Range("I1").Value = "id de la rémunération"

# Row 2: A 361 -> 362, C 16 -> 13, new H2 = 0.4375, new I2 = 43
$wsCollab.Range("A2").Value = 362
$wsCollab.Range("C2").Value = 13
$wsCollab.Range("H2").Value = 0.4375
$wsCollab.Range("I2").Value = 43

# New row 3
$wsCollab.Range("A3").Value = 363
$wsCollab.Range("B3").Value = 574
$wsCollab.Range("C3").Value = 13
$wsCollab.Range("D3").Value = 7
$wsCollab.Range("E3").Value = 4
$wsCollab.Range("F3").Value = 1
$wsCollab.Range("G3").Value = 0.4236111111111111
$wsCollab.Range("G3").NumberFormat = "h:mm:ss"
$wsCollab.Range("H3").Value = 0.43055555555555558
$wsCollab.Range("H3").NumberFormat = "h:mm:ss"
$wsCollab.Range("I3").Value = 44

# New row 4
$wsCollab.Range("A4").Value = 364
$wsCollab.Range("B4").Value = 574
$wsCollab.Range("C4").Value = 13
$wsCollab.Range("D4").Value = 7
$wsCollab.Range("E4").Value = 4
$wsCollab.Range("F4").Value = 1
$wsCollab.Range("G4").Value = 0.98611111111111116
$wsCollab.Range("G4").NumberFormat = "h:mm:ss"
$wsCollab.Range("H4").Value = 0.069444444444444434
$wsCollab.Range("H4").NumberFormat = "h:mm:ss"
$wsCollab.Range("I4").Value = 45

# Widen the new column (closest achievable width to 20.7109375 chars)
$wsCollab.Columns.Item(9).ColumnWidth = 19.8

# Restore Collab as the active / displayed sheet (it was the tab shown
# before the edits, and selecting on Chauffeur above would otherwise
# switch the active tab to it).
$wsCollab.Activate()

Write-Output "edits applied"
